$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-CellText($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $find = $cell.Range.Find
    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# Row 2 (Medicações controladas -> Tempo restrito...)
Replace-CellText 2 1 "Medicações controladas com receita especial" "Tempo restrito para entrega do sistema"
Replace-CellText 2 2 "De acordo com a legislação do MAPA, para algumas medicações, há a necessidade de receita especial como alguns analgésicos ou sedativos. Mesmo os atuais E-commerces para produtos pets não disponibilizam essas medicações para compra online." "O sistema deve ficar pronto no período de um ano."

# Row 3 (Vigilância sanitária -> Sem orçamento...)
Replace-CellText 3 1 "Vigilância sanitária" "Sem orçamento para o desenvolvimento do sistema."
Replace-CellText 3 2 "A loja deve ficar muito atenta na venda de produtos que possuam data de validade como rações, petiscos, suplementos e medicamentos. O envio de um produto com a data de validade vencida pode levar a processos judiciais a loja e perdas financeiras." "Não há orçamento para o desenvolvimento do sistema. A equipe terá que utilizar ferramentas open source para o desenvolvimento do sistema. Não haverá a possibilidade do uso de ferramentas pagas."

# Row 4 (A linguagem escolhida...TypeScript. -> Não poder se comunicar...)
Replace-CellText 4 1 "A linguagem escolhida para o desenvolvimento será TypeScript." "Não poder se comunicar com cliente em alguns dias."
Replace-CellText 4 2 "Com foco em performance e produtividade decidimos realizar está escolha em conjunto com o cliente visando a gama de benefícios que nos oferece." "O cliente pediu para que em alguns dias a equipe não poderá se comunicar com o cliente."

# Remove the last two rows (banco de dados / Web rows) entirely.
# Delete from the end so earlier indices stay valid.
$t.Rows.Item($t.Rows.Count).Delete()
$t.Rows.Item($t.Rows.Count).Delete()

Write-Output "done"
